$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.184474945068359
$ws.Range("B1").Value = 2.410895109176636
$ws.Range("C1").Value = 1.43298065662384
$ws.Range("D1").Value = 1.506852746009827
$ws.Range("E1").Value = 1.597444176673889
